# Apply "DeudoresPrueba" update:
#  - Add new client ALAMO as a new first data row (list is kept alphabetical by client)
#  - Remove the two CANTON WOK rows
#  - Remove the MERKA FRUVER ALEJANDRO row
#  - Add new client COCINA CHINA (alphabetical order, after CLIENTE PAOLA)
#  - Remove the PLANADAS NUEVO row
#  - Update PUNTA DE ANCA's date and value
#
# Final data block (rows 2-31) written explicitly, then leftover trailing
# rows from the longer original table are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{A=1; B='ALAMO'; C=46027; D=288000},
    @{A=2; B='ALISO'; C=46022; D=93000},
    @{A=3; B='ARROZ PAISA SUBA'; C=46022; D=166000},
    @{A=4; B='CAMPO VERDE TOCANCIPA'; C=46021; D=475000},
    @{A=5; B='CAMPO VERDE ZIPAQUIRA'; C=46021; D=18900},
    @{A=8; B='CARNES JOHANA'; C=46021; D=320000},
    @{A=9; B='CIMARRON DORADO'; C=46020; D=449800},
    @{A=10; B='CLIENTE PAOLA'; C=46018; D=174000},
    @{A=11; B='COCINA CHINA'; C=46027; D=170000},
    @{A=12; B='CRISTIAN ACACIAS'; C=46009; D=1000000},
    @{A=13; B='DARWIN FUTBOL'; C=45921; D=200000},
    @{A=14; B='DAVIDCITO'; C=45947; D=100000},
    @{A=15; B='EL JORDAN'; C=46022; D=1600000},
    @{A=16; B='FRANCO'; C=45996; D=20000},
    @{A=17; B='FRANCO'; C=46017; D=545800},
    @{A=18; B='LA CABAÑA'; C=46020; D=215300},
    @{A=19; B='LA PAMPA'; C=46006; D=229900},
    @{A=20; B='LA SELECTA'; C=45912; D=82000},
    @{A=21; B='MAFE'; C=46017; D=190000},
    @{A=22; B='MERKA FRUVER DEXI'; C=45988; D=15400},
    @{A=23; B='MERKA FRUVER DEXI'; C=45995; D=339000},
    @{A=24; B='MICHAEL'; C=46011; D=80000},
    @{A=25; B='NEVADA'; C=46020; D=195000},
    @{A=26; B='PARAÍSO FUNZA'; C=46020; D=276000},
    @{A=27; B='PARAÍSO MOSQUERA'; C=46013; D=328800},
    @{A=28; B='PINILLA'; C=45931; D=82000},
    @{A=29; B='PUNTA DE ANCA'; C=46027; D=200000},
    @{A=30; B='SAMY 2'; C=46021; D=203000},
    @{A=31; B='SAMY 2'; C=46013; D=142000},
    @{A=32; B='WILINTONG'; C=46006; D=100000}
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $false
    $r++
}

# The original table had 32 data rows (through row 33); the updated table
# only has 30 (through row 31) -- clear the two now-unused trailing rows.
$ws.Range("A32:E33").Clear()
